# Add the new result row (num classes = 1? actually A4 = 1 per diff, the "10 classes" refers to
# the new podcast list having 10 shows) with the wider podcast-name list and new accuracy value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "2 bears one cave, the views, andrew schulz, impaulsive, joe rogan, h3h3, jenna and julien, lex fridman, off topic, optic podcast"
$ws.Range("C4").Value = 96.5

# Widen column B to fit the longer podcast-name string (matches author's manual resize).
$ws.Columns.Item(2).ColumnWidth = 103.89322916666667

# Matches the saved selection state in the edited workbook.
$ws.Range("C6").Select()
